$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "budget_group" between "budget" and "kpi"
# ------------------------------------------------------------------
$budgetSheet = $wb.Worksheets.Item("budget")

$newSheet = $wb.Worksheets.Add()
$newSheet.Move($null, $budgetSheet) | Out-Null
$newSheet.Name = "budget_group"

# Re-fetch the "kpi" sheet reference by name now that the sheet
# collection has been reordered (index-based references captured
# before the insert would now point at the wrong sheet).
$kpiSheet = $wb.Worksheets.Item("kpi")

# ------------------------------------------------------------------
# 2. Populate the new "budget_group" sheet
# ------------------------------------------------------------------
$newSheet.Range("A1").Value = "section"
$newSheet.Range("B1").Value = "project"
$newSheet.Range("C1").Value = "budget"
$newSheet.Range("D1").Value = "actual"

$newSheet.Range("A2").Value = "Alpha"
$newSheet.Range("B2").Value = "A"
$newSheet.Range("C2").Value = 100
$newSheet.Range("D2").Value = 75

$newSheet.Range("A3").Value = "Alpha"
$newSheet.Range("B3").Value = "B"
$newSheet.Range("C3").Value = 120
$newSheet.Range("D3").Value = 60

$newSheet.Range("A4").Value = "Alpha"
$newSheet.Range("B4").Value = "C"
$newSheet.Range("C4").Value = 80
$newSheet.Range("D4").Value = 80

$newSheet.Range("A5").Value = "Beta"
$newSheet.Range("B5").Value = "D"
$newSheet.Range("C5").Value = 55
$newSheet.Range("D5").Value = 65

$newSheet.Range("A6").Value = "Beta"
$newSheet.Range("B6").Value = "E"
$newSheet.Range("C6").Value = 35
$newSheet.Range("D6").Value = 70

$newSheet.Range("A7").Value = "Beta"
$newSheet.Range("B7").Value = "F"
$newSheet.Range("C7").Value = 75
$newSheet.Range("D7").Value = 75

$newSheet.Range("A8").Value = "Beta"
$newSheet.Range("B8").Value = "G"
$newSheet.Range("C8").Value = 90
$newSheet.Range("D8").Value = 45

$newSheet.Range("A9").Value = "Gamma"
$newSheet.Range("B9").Value = "H"
$newSheet.Range("C9").Value = 65
$newSheet.Range("D9").Value = 65

$newSheet.Range("A10").Value = "Gamma"
$newSheet.Range("B10").Value = "I"
$newSheet.Range("C10").Value = 55
$newSheet.Range("D10").Value = 60

$newSheet.Range("A1:D10").Columns.AutoFit()

# Match page / print setup of the sibling "budget" sheet.
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1
$newSheet.PageSetup.RightHeader = "&P/&N"
$newSheet.PageSetup.LeftFooter = "&8&Z&F"
$newSheet.PageSetup.RightFooter = "&8&D &T"

# selection / active cell on the new sheet
$newSheet.Range("C5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130

# ------------------------------------------------------------------
# 3. Update data on the "budget" sheet
# ------------------------------------------------------------------
$budgetSheet.Range("C4").Value = 80
$budgetSheet.Range("B5").Value = 55
$budgetSheet.Range("C5").Value = 65

$budgetSheet.Range("A2:C5").Select() | Out-Null

# ------------------------------------------------------------------
# 4. Update selection on the "kpi" sheet (kept as-is: B3)
# ------------------------------------------------------------------
$kpiSheet.Range("B3").Select() | Out-Null

# ------------------------------------------------------------------
# 5. Activate the new "budget_group" sheet (becomes the active tab)
# ------------------------------------------------------------------
$newSheet.Activate() | Out-Null
